$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.370.20"
$ws.Range("E2").Value = "  -2.89%  "

# Row 3
$ws.Range("D3").Value = "3.090.01"
$ws.Range("E3").Value = "  -1.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "548.67"
$ws.Range("E5").Value = "  -3.09%  "

# Row 6
$ws.Range("D6").Value = "137.42"
$ws.Range("E6").Value = "  -7.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "3.082.28"
$ws.Range("E8").Value = "  -1.63%  "

# Row 9
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -0.90%  "

# Row 10
$ws.Range("D10").Value = "6.63"
$ws.Range("E10").Value = "  -3.99%  "

# Row 11
$ws.Range("E11").Value = "  -0.25%  "

# Row 12
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$ws.Range("D13").Value = "35.08"
$ws.Range("E13").Value = "  -3.03%  "

# Row 14
$ws.Range("D14").Value = "0.0000217"
$ws.Range("E14").Value = "  -1.93%  "

# Row 15
$ws.Range("D15").Value = "3.590.38"
$ws.Range("E15").Value = "  -1.50%  "

# Row 16
$ws.Range("D16").Value = "63.408.22"
$ws.Range("E16").Value = "  -2.87%  "

# Row 17
$ws.Range("E17").Value = "  -1.01%  "

# Row 18
$ws.Range("D18").Value = "3.094.12"
$ws.Range("E18").Value = "  -1.60%  "

# Row 19
$ws.Range("D19").Value = "6.67"
$ws.Range("E19").Value = "  -1.04%  "

# Row 20
$ws.Range("D20").Value = "485.66"
$ws.Range("E20").Value = "  -7.44%  "

# Row 21
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -2.97%  "

# Row 22
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "7.19"
$ws.Range("E23").Value = "  -2.76%  "

# Row 24
$ws.Range("D24").Value = "77.83"
$ws.Range("E24").Value = "  -1.12%  "

# Row 25
$ws.Range("D25").Value = "12.25"
$ws.Range("E25").Value = "  -3.78%  "

# Row 26
$ws.Range("E26").Value = "  +0.29%  "

# Row 27
$ws.Range("E27").Value = "  -2.29%  "

# Row 28
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  -4.02%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  -8.79%  "

# Row 31
$ws.Range("D31").Value = "26.46"
$ws.Range("E31").Value = "  +1.07%  "

# Row 32
$ws.Range("E32").Value = "  -0.64%  "

# Row 33
$ws.Range("D33").Value = "2.51"
$ws.Range("E33").Value = "  -6.39%  "

# Row 34
$ws.Range("D34").Value = "60.68"
$ws.Range("E34").Value = "  +14.41%  "

# Row 35
$ws.Range("D35").Value = "533.91"
$ws.Range("E35").Value = "  -4.17%  "

# Row 36
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -1.79%  "

# Row 37
$ws.Range("D37").Value = "5.18"
$ws.Range("E37").Value = "  -5.11%  "

# Row 38
$ws.Range("D38").Value = "0.0402"
$ws.Range("E38").Value = "  -9.12%  "

# Row 39
$ws.Range("D39").Value = "0.0794"
$ws.Range("E39").Value = "  -3.45%  "

# Row 40
$ws.Range("D40").Value = "3.067.84"
$ws.Range("E40").Value = "  -0.31%  "

# Row 41
$ws.Range("E41").Value = "  -2.10%  "

# Row 42
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  -7.43%  "

# Row 43
$ws.Range("D43").Value = "8.12"
$ws.Range("E43").Value = "  -1.57%  "

# Row 44
$ws.Range("D44").Value = "0.255"
$ws.Range("E44").Value = "  -1.09%  "

# Row 46
$ws.Range("D46").Value = "2.05"
$ws.Range("E46").Value = "  -6.71%  "

# Row 47
$ws.Range("D47").Value = "121.52"
$ws.Range("E47").Value = "  +2.90%  "

# Row 48
$ws.Range("D48").Value = "24.38"
$ws.Range("E48").Value = "  -2.58%  "

# Row 49
$ws.Range("E49").Value = "  -1.39%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0506"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +61.14%  "
